$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 88
$ws.Range("F4").Value = 632
$ws.Range("F5").Value = 355
$ws.Range("F8").Value = 11585
$ws.Range("F9").Value = 195
$ws.Range("F12").Value = 2090
$ws.Range("F16").Value = 232
$ws.Range("F18").Value = 1184
$ws.Range("F19").Value = 160
$ws.Range("F20").Value = 242
$ws.Range("F21").Value = 731
$ws.Range("F23").Value = 260
$ws.Range("F25").Value = 718
$ws.Range("F26").Value = 3478
$ws.Range("F27").Value = 1052
$ws.Range("F28").Value = 783
$ws.Range("F32").Value = 973
$ws.Range("F35").Value = 253
$ws.Range("F36").Value = 16
$ws.Range("F38").Value = 4
$ws.Range("F39").Value = 2336
$ws.Range("F40").Value = 4379
$ws.Range("F41").Value = 5449
$ws.Range("F42").Value = 105
$ws.Range("F45").Value = 254
$ws.Range("F46").Value = 60
$ws.Range("F47").Value = 25
$ws.Range("F49").Value = 94

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F11").Value = 606
$ws.Range("F23").Value = 5

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 749
$ws.Range("F3").Value = 417
$ws.Range("F4").Value = 59

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 749
$ws.Range("F3").Value = 417
$ws.Range("F4").Value = 59
$ws.Range("F5").Value = 88
$ws.Range("F7").Value = 632
$ws.Range("F8").Value = 355
$ws.Range("F10").Value = 11585
$ws.Range("F13").Value = 2090
$ws.Range("F18").Value = 1184
$ws.Range("F19").Value = 160
$ws.Range("F20").Value = 242
$ws.Range("F23").Value = 731
$ws.Range("F24").Value = 718
$ws.Range("F26").Value = 783
$ws.Range("F31").Value = 973
$ws.Range("F33").Value = 253
$ws.Range("F34").Value = 16
$ws.Range("F36").Value = 4379
$ws.Range("F38").Value = 105
$ws.Range("F41").Value = 254
$ws.Range("F44").Value = 60
$ws.Range("F48").Value = 94
$ws.Range("F49").Value = 5

